# Auto-generated Excel COM-interop script
# Updates the cryptocurrency price/volume table to reflect the latest
# scrape from coinranking.com (GitHub Actions scheduled refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.419.50"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "2.911.43"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'352.83"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "'112.03"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.561"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.628"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'39.95"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").Value = "'0.0867"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("D12").Value = "'0.136"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'19.86"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "'7.79"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "3.367.13"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "'1.01"
$ws.Range("E16").Value = "  +7.09%  "
$ws.Range("D17").Value = "2.904.23"
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("D18").Value = "52.385.81"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "'7.63"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "'3.32"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "'14.21"
$ws.Range("E21").Value = "  +4.39%  "
$ws.Range("D22").Value = "0.0₃0980"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'70.81"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "'269.80"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'2.77"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "'26.74"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.168"
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'10.65"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'6.35"
$ws.Range("E30").Value = "  +12.41%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'37.79"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").Value = "'6.61"
$ws.Range("E32").Value = "  +7.64%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0989"
$ws.Range("E33").Value = "  +11.85%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.25"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'53.41"
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("D36").Value = "'0.0451"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +5.47%  "
$ws.Range("D39").Value = "'18.81"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "  +13.72%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.117"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'23.53"
$ws.Range("E43").Value = "  +7.07%  "
$ws.Range("D44").Value = "'2.66"
$ws.Range("E44").Value = "  +10.59%  "
$ws.Range("D45").Value = "'120.76"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'3.55"
$ws.Range("E47").Value = "  +4.12%  "
$ws.Range("D48").Value = "2.199.05"
$ws.Range("E48").Value = "  +4.36%  "
$ws.Range("D49").Value = "'0.267"
$ws.Range("E49").Value = "  +24.68%  "
$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").Value = "'0.970"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "'0.0336"
$ws.Range("E51").Value = "  +11.17%  "
